# Sheet "汽車" (Car) is the 3rd worksheet in this workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# --- Row 1 currently duplicates row 2's data; turn it into a proper header
#     row (matching the header convention used on every other property sheet:
#     name, capacity, owner, register_date, register_reason, acquire_value,
#     property_category, category, date, legislator_name, legislator_id,
#     source_file, index). ---
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# Give the new header cells (H1:N1) the same bold / bordered / centered
# look already used by the existing header cells B1:G1.
$hdr = $ws.Range("H1:N1")
$hdr.Font.Bold = $true
$hdr.Borders.LineStyle = 1
$hdr.HorizontalAlignment = -4108   # xlCenter
$hdr.VerticalAlignment = -4160     # xlTop

# --- Row 2: TOYOTACAMRYLE (index 47) ---
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"
# "date" column holds a literal text date string, not a real Excel date -
# force text format first so Excel doesn't auto-convert it to a serial date.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "2012-04-23"
$ws.Range("K2").Value = "許忠信"
$ws.Range("L2").Value = 1749
$ws.Range("M2").Value = "tmp50641"
$ws.Range("N2").Value = 47

# --- Row 3: 中華GL20SS58 (index 48) ---
$ws.Range("H3").Value = "land"
$ws.Range("I3").Value = "normal"
$ws.Range("J3").NumberFormat = "@"
$ws.Range("J3").Value = "2012-04-23"
$ws.Range("K3").Value = "許忠信"
$ws.Range("L3").Value = 1749
$ws.Range("M3").Value = "tmp50641"
$ws.Range("N3").Value = 48
